# "Generate Report for Handoff"
# Updates the localization-status report: flips the Status column from
# "Handed back: in sync with en-US" to "Ready for handoff", refreshes the
# "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" timestamps,
# and narrows the (now shorter) Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet: zh-cn / de-de status columns + generate-date column ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-19 17:03:04"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-19 17:02:56"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-19 17:03:04"

# --- Narrow the Status columns now that the text is shorter ---
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
